$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 68

$ws.Cells.Item($newRow, 1).NumberFormat = "@"
$ws.Cells.Item($newRow, 1).Value = "2020-08-06"
$ws.Cells.Item($newRow, 1).Style = "Normal"
$ws.Cells.Item($newRow, 2).Value = 462690
$ws.Cells.Item($newRow, 3).Value = 506252
$ws.Cells.Item($newRow, 4).Value = 87973
$ws.Cells.Item($newRow, 5).Value = 50517
$ws.Cells.Item($newRow, 6).Value = 26.8
